# Fruta / hortaliza, semanal
# Insert a new weekly observation row before row 48 (shifting the existing
# rows 48-106 down to 49-107) and populate the new row with the latest
# data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 48..106 down to 49..107
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new weekly record
$ws.Range("A48").Value = 9
$ws.Range("B48").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C48").Value = "Metropolitana"
$ws.Range("D48").Value = 44482
$ws.Range("E48").Value = 13
$ws.Range("F48").Value = "Fruta"
$ws.Range("G48").Value = 100101
$ws.Range("H48").Value = "Berries"
$ws.Range("I48").Value = 100101001
$ws.Range("J48").Value = "Arándano (blue)"
$ws.Range("K48").Value = "Sin especificar"
$ws.Range("L48").Value = "Primera"
$ws.Range("M48").Value = 220
$ws.Range("N48").Value = 14000
$ws.Range("O48").Value = 14000
$ws.Range("P48").Value = 14000
$ws.Range("Q48").Value = "$/bandeja 2 kilos"
$ws.Range("R48").Value = "Provincia del Elquí"
$ws.Range("S48").Value = 7000
$ws.Range("T48").Value = 2
